$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 170-171),
# pushing the existing rows 170-182 down to 172-184.
$ws.Rows("170:171").Insert()

# Row 170: new weekly price entry ($/caja 18 kilos)
$ws.Range("A170").Value = 9
$ws.Range("B170").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C170").Value = "Metropolitana"
$ws.Range("D170").Value = 45166
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = 100114002
$ws.Range("G170").Value = "Camote"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 700
$ws.Range("K170").Value = 18000
$ws.Range("L170").Value = 19000
$ws.Range("M170").Value = 18500
$ws.Range("N170").Value = "$/caja 18 kilos"
$ws.Range("O170").Value = "Perú"
$ws.Range("P170").Value = 1028
$ws.Range("Q170").Value = 18
$ws.Range("R170").Value = "Hortaliza"

# Row 171: new weekly price entry ($/malla 18 kilos)
$ws.Range("A171").Value = 9
$ws.Range("B171").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C171").Value = "Metropolitana"
$ws.Range("D171").Value = 45166
$ws.Range("E171").Value = 13
$ws.Range("F171").Value = 100114002
$ws.Range("G171").Value = "Camote"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 610
$ws.Range("K171").Value = 15000
$ws.Range("L171").Value = 16000
$ws.Range("M171").Value = 15500
$ws.Range("N171").Value = "$/malla 18 kilos"
$ws.Range("O171").Value = "Perú"
$ws.Range("P171").Value = 861
$ws.Range("Q171").Value = 18
$ws.Range("R171").Value = "Hortaliza"
